# chore: consider header for column-wise drop
#
# Insert a new row above the existing row 4 (pushing the former rows
# 4-24 down to 5-25) and seed the new A4 with a =NA() formula so the
# sheet now has an explicit "#N/A" marker row right under the header,
# instead of the first data row being silently treated as data for a
# column-wise drop.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 4:24 down to 5:25, leaving a blank row 4 behind.
$ws.Rows.Item(4).Insert()

# New row 4 only gets a formula in column A.
$ws.Range("A4").Formula = "=NA()"

# Match the author's recorded selection after making the edit.
$ws.Range("B4").Select()
